$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1082.0667
$ws.Range("I2").Value = 282.55554
$ws.Range("J2").Value = 2281.3333
$ws.Range("K2").Value = 282.55554
$ws.Range("L2").Value = 2281.3333
$ws.Range("M2").Value = -169.55554
$ws.Range("N2").Value = -2507.3333
$ws.Range("H5").Value = 105
$ws.Range("I5").Value = 105
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 105
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 10
$ws.Range("N5").ClearContents()
$ws.Range("H41").Value = 440.42856
$ws.Range("I41").Value = 436.6
$ws.Range("K41").Value = 436.6
$ws.Range("M41").Value = 3.399999999999977
$ws.Range("H64").Value = 24103
$ws.Range("I64").Value = 3890
$ws.Range("K64").Value = 3890
$ws.Range("M64").Value = -3642
$ws.Range("H67").Value = 24103
$ws.Range("I67").Value = 3890
$ws.Range("K67").Value = 3890
$ws.Range("M67").Value = -3032
$ws.Range("H74").Value = 4398.8
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 4997
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 4997
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -6869
$ws.Range("H77").Value = 4398.8
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 4997
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 24985
$ws.Range("M77").Value = -15320
$ws.Range("N77").Value = -34345
$ws.Range("H113").Value = 3690.4546
$ws.Range("I113").Value = 3033
$ws.Range("K113").Value = 3033
$ws.Range("M113").Value = 221
$ws.Range("H137").Value = 4979.6
$ws.Range("I137").Value = 934.2222
$ws.Range("J137").Value = 11047.667
$ws.Range("K137").Value = 2802.6666
$ws.Range("L137").Value = 33143.001
$ws.Range("M137").Value = -252.6666
$ws.Range("N137").Value = -38243.001
$ws.Range("H138").Value = 3504.525
$ws.Range("I138").Value = 4040.2727
$ws.Range("J138").Value = 3301.3103
$ws.Range("K138").Value = 12120.8181
$ws.Range("L138").Value = 9903.930899999999
$ws.Range("M138").Value = -6980.8181
$ws.Range("N138").Value = -20183.9309
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 54926.137
$ws.Range("I32").Value = 70314.375
$ws.Range("J32").Value = 13890.833
$ws.Range("K32").Value = 70314.375
$ws.Range("L32").Value = 13890.833
$ws.Range("M32").Value = -70027.375
$ws.Range("N32").Value = -14464.833
$ws.Range("H61").Value = 4446.1665
$ws.Range("I61").Value = 4414
$ws.Range("K61").Value = 4414
$ws.Range("M61").Value = -4202
$ws.Range("H102").Value = 5125
$ws.Range("I102").Value = 3500
$ws.Range("K102").Value = 3500
$ws.Range("M102").Value = -1878
$ws.Range("H110").Value = 2095
$ws.Range("I110").Value = 2073.158
$ws.Range("J110").Value = 2233.3333
$ws.Range("K110").Value = 2073.158
$ws.Range("L110").Value = 2233.3333
$ws.Range("M110").Value = -28.1579999999999
$ws.Range("N110").Value = -6323.3333
$ws.Range("H111").Value = 27500
$ws.Range("J111").Value = 27500
$ws.Range("L111").Value = 27500
$ws.Range("N111").Value = -35680
$ws.Range("H132").Value = 171614.17
$ws.Range("J132").Value = 3999
$ws.Range("L132").Value = 11997
$ws.Range("N132").Value = -17057
$ws.Range("H136").Value = 4446.1665
$ws.Range("I136").Value = 4414
$ws.Range("K136").Value = 13242
$ws.Range("M136").Value = -10692
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4973.6
$ws.Range("I94").Value = 5827.8184
$ws.Range("J94").Value = 2624.5
$ws.Range("K94").Value = 5827.8184
$ws.Range("L94").Value = 2624.5
$ws.Range("M94").Value = -5376.8184
$ws.Range("N94").Value = -3526.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 14.2
$ws.Range("I25").Value = 14.2
$ws.Range("K25").Value = 14.2
$ws.Range("M25").Value = 159.8
$ws.Range("H58").Value = 74546.64
$ws.Range("I58").Value = 74546.64
$ws.Range("K58").Value = 74546.64
$ws.Range("M58").Value = -74343.64
$ws.Range("H127").Value = 80900
$ws.Range("J127").Value = 80900
$ws.Range("L127").Value = 80900
$ws.Range("N127").Value = -90820
$ws.Range("H132").Value = 1254.8948
$ws.Range("I132").Value = 1155.8909
$ws.Range("J132").Value = 3977.5
$ws.Range("K132").Value = 3467.6727
$ws.Range("L132").Value = 11932.5
$ws.Range("M132").Value = -937.6727000000001
$ws.Range("N132").Value = -16992.5
$ws.Range("H136").Value = 74546.64
$ws.Range("I136").Value = 74546.64
$ws.Range("K136").Value = 223639.92
$ws.Range("M136").Value = -221089.92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 559113.4
$ws.Range("I4").Value = 859574.4
$ws.Range("K4").Value = 2578723.2
$ws.Range("M4").Value = -2578611.2
$ws.Range("H29").Value = 645.2286
$ws.Range("I29").Value = 570.9474
$ws.Range("J29").Value = 733.4375
$ws.Range("K29").Value = 1712.8422
$ws.Range("L29").Value = 2200.3125
$ws.Range("M29").Value = -1435.8422
$ws.Range("N29").Value = -2754.3125
$ws.Range("H41").Value = 470
$ws.Range("I41").Value = 450
$ws.Range("J41").Value = 550
$ws.Range("K41").Value = 1350
$ws.Range("L41").Value = 1650
$ws.Range("M41").Value = -1012
$ws.Range("N41").Value = -2326
$ws.Range("H69").Value = 973.3333
$ws.Range("I69").Value = 973.25
$ws.Range("J69").Value = 974
$ws.Range("K69").Value = 2919.75
$ws.Range("L69").Value = 2922
$ws.Range("M69").Value = -2108.75
$ws.Range("N69").Value = -4544
$ws.Range("H72").Value = 973.3333
$ws.Range("I72").Value = 973.25
$ws.Range("J72").Value = 974
$ws.Range("K72").Value = 8759.25
$ws.Range("L72").Value = 8766
$ws.Range("M72").Value = -4703.25
$ws.Range("N72").Value = -16878
$ws.Range("H92").Value = 457.96774
$ws.Range("I92").Value = 371.7143
$ws.Range("J92").Value = 529
$ws.Range("K92").Value = 1115.1429
$ws.Range("L92").Value = 1587
$ws.Range("M92").Value = 132.8571000000002
$ws.Range("N92").Value = -4083
$ws.Range("H122").Value = 527.4
$ws.Range("I122").Value = 434.25
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 3908.25
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -1458.25
$ws.Range("N122").Value = -13000
$ws.Range("H131").Value = 3854471
$ws.Range("I131").Value = 2472.5557
$ws.Range("J131").Value = 5893764
$ws.Range("K131").Value = 7417.6671
$ws.Range("L131").Value = 17681292
$ws.Range("M131").Value = -2377.6671
$ws.Range("N131").Value = -17691372
$ws.Range("H132").Value = 1474
$ws.Range("J132").Value = 1315.3334
$ws.Range("L132").Value = 11838.0006
$ws.Range("N132").Value = -16898.0006
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8559.6
$ws.Range("J70").Value = 6596.3335
$ws.Range("L70").Value = 6596.3335
$ws.Range("N70").Value = -7136.3335
$ws.Range("H73").Value = 8559.6
$ws.Range("J73").Value = 6596.3335
$ws.Range("L73").Value = 6596.3335
$ws.Range("N73").Value = -8468.333500000001
$ws.Range("H97").Value = 1590.1111
$ws.Range("I97").Value = 1590.1111
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1590.1111
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1094.1111
$ws.Range("N97").ClearContents()
$ws.Range("H101").Value = 57499.5
$ws.Range("J101").Value = 57499.5
$ws.Range("L101").Value = 57499.5
$ws.Range("N101").Value = -63989.5
$ws.Range("H122").Value = 2116.6128
$ws.Range("I122").Value = 2084.92
$ws.Range("J122").Value = 2248.6667
$ws.Range("K122").Value = 6254.76
$ws.Range("L122").Value = 6746.000100000001
$ws.Range("M122").Value = -3804.76
$ws.Range("N122").Value = -11646.0001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3898
$ws.Range("I7").Value = 3851.1538
$ws.Range("K7").Value = 3851.1538
$ws.Range("M7").Value = -3739.1538
$ws.Range("H122").Value = 3836.4644
$ws.Range("I122").Value = 3339.2222
$ws.Range("J122").Value = 4072
$ws.Range("K122").Value = 10017.6666
$ws.Range("L122").Value = 12216
$ws.Range("M122").Value = -7567.6666
$ws.Range("N122").Value = -17116
$ws.Range("H126").Value = 3898
$ws.Range("I126").Value = 3851.1538
$ws.Range("K126").Value = 11553.4614
$ws.Range("M126").Value = -9083.4614
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2497.875
$ws.Range("I96").Value = 2497.875
$ws.Range("K96").Value = 2497.875
$ws.Range("M96").Value = -1124.875
$ws.Range("H122").Value = 1286
$ws.Range("I122").Value = 1337.1818
$ws.Range("J122").Value = 1004.5
$ws.Range("K122").Value = 4011.5454
$ws.Range("L122").Value = 3013.5
$ws.Range("M122").Value = -1561.5454
$ws.Range("N122").Value = -7913.5
$ws.Range("H132").Value = 75322.28
$ws.Range("I132").Value = 83373.28999999999
$ws.Range("J132").Value = 4876
$ws.Range("K132").Value = 250119.87
$ws.Range("L132").Value = 14628
$ws.Range("M132").Value = -247589.87
$ws.Range("N132").Value = -19688
$ws.Range("H136").Value = 2142.025
$ws.Range("I136").Value = 1825.125
$ws.Range("K136").Value = 5475.375
$ws.Range("M136").Value = -2925.375
